# Update the cryptos list with freshly scraped Price / Volume(1h) values.
# Columns: A=index, B=Coin, C=Link, D=Price, E=Volume(1h)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new Price, new Volume(1h)). $null means the Price cell is
# unchanged for that row (only Volume(1h) was updated).
$updates = @{
    2  = @{ D = "27.555.95";  E = "  -0.17%  " }
    3  = @{ D = "1.583.51";   E = "  -0.88%  " }
    4  = @{ D = $null;        E = "  -0.20%  " }
    5  = @{ D = "208.34";     E = "  -0.14%  " }
    6  = @{ D = "0.498";      E = "  -1.04%  " }
    7  = @{ D = $null;        E = "  -0.14%  " }
    8  = @{ D = $null;        E = "  +0.12%  " }
    10 = @{ D = $null;        E = "  -0.35%  " }
    11 = @{ D = "0.0866";     E = "  +0.20%  " }
    12 = @{ D = "1.810.05";   E = "  -0.84%  " }
    13 = @{ D = "1.570.61";   E = "  -1.59%  " }
    14 = @{ D = "3.83";       E = "  -0.99%  " }
    15 = @{ D = "0.525";      E = "  -2.42%  " }
    16 = @{ D = "27.546.27";  E = "  -0.23%  " }
    17 = @{ D = "63.10";      E = "  -0.69%  " }
    18 = @{ D = "215.82";     E = "  -1.00%  " }
    19 = @{ D = "7.34";       E = "  -1.04%  " }
    20 = @{ D = $null;        E = "  -0.61%  " }
    21 = @{ D = "1.00";       E = "  -0.18%  " }
    22 = @{ D = $null;        E = "  -1.23%  " }
    23 = @{ D = $null;        E = "  +1.17%  " }
    24 = @{ D = "2.02";       E = "  +0.07%  " }
    25 = @{ D = "153.30";     E = "  -0.16%  " }
    26 = @{ D = "6.92";       E = "  +2.61%  " }
    27 = @{ D = $null;        E = "  -0.16%  " }
    28 = @{ D = "15.08";      E = "  -0.14%  " }
    30 = @{ D = $null;        E = "  -0.60%  " }
    31 = @{ D = $null;        E = "  +0.69%  " }
    32 = @{ D = "3.23";       E = "  -1.23%  " }
    33 = @{ D = "1.377.28";   E = "  +0.30%  " }
    34 = @{ D = $null;        E = "  -0.27%  " }
    35 = @{ D = $null;        E = "  +0.71%  " }
    36 = @{ D = "0.973";      E = "  +0.14%  " }
    37 = @{ D = $null;        E = "  +0.01%  " }
    38 = @{ D = $null;        E = "  +1.23%  " }
    39 = @{ D = "0.535";      E = "  -1.15%  " }
    40 = @{ D = "0.828";      E = "  +1.49%  " }
    41 = @{ D = $null;        E = "  -0.16%  " }
    42 = @{ D = "0.970";      E = "  -0.62%  " }
    43 = @{ D = "1.80";       E = "  +0.96%  " }
    44 = @{ D = "64.55";      E = "  +0.60%  " }
    45 = @{ D = "5.30";       E = "  -1.54%  " }
    46 = @{ D = $null;        E = "  +3.06%  " }
    47 = @{ D = "1.720.27";   E = "  -0.90%  " }
    48 = @{ D = "85.42";      E = "  -2.87%  " }
    49 = @{ D = "0.0₇0994";   E = "  -0.81%  " }
    50 = @{ D = "0.0960";     E = "  -1.13%  " }
    51 = @{ D = "0.0495";     E = "  -0.66%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals.D) {
        # Force text storage so values like "1.00" / "63.10" keep their
        # original formatting instead of being coerced to numbers (1 / 63.1).
        $ws.Range("D$row").NumberFormat = "@"
        $ws.Range("D$row").Value = $vals.D
    }
    $ws.Range("E$row").Value = $vals.E
}
